$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 45204 to 45205 for every data row (2..518)
for ($r = 2; $r -le 518; $r++) {
    $ws.Cells.Item($r, 3).Value = 45205
}

# Row 518 picks up an explicit row height (15, customHeight) in the new file
$ws.Rows.Item(518).RowHeight = 15

# New row 519 appended with the latest filing
$ws.Range("A519").Value = "A 47727-2023"

$ws.Range("B519").Value = 45203
$ws.Range("B519").NumberFormat = "YYYY-MM-DD"

$ws.Range("C519").Value = 45205
$ws.Range("C519").NumberFormat = "YYYY-MM-DD"

$ws.Range("D519").Value = "JÄMTLANDS LÄN"
$ws.Range("E519").Value = "ÖSTERSUND"

$ws.Range("G519").Value = 2.6
$ws.Range("H519").Value = 0
$ws.Range("I519").Value = 0
$ws.Range("J519").Value = 0
$ws.Range("K519").Value = 0
$ws.Range("L519").Value = 0
$ws.Range("M519").Value = 0
$ws.Range("N519").Value = 0
$ws.Range("O519").Value = 0
$ws.Range("P519").Value = 0
$ws.Range("Q519").Value = 0

$ws.Range("R519").Value = ""
$ws.Range("R519").WrapText = $true
